# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp label in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 10:35"

# Row 34 (Polonia)
$ws.Range("B34").Value = 17469
$ws.Range("C34").Value = 265
$ws.Range("E34").Value = 9904
$ws.Range("G34").Value = 8
$ws.Range("H34").Value = 869

# Row 43 (Filipinas)
$ws.Range("B43").Value = 11876
$ws.Range("C43").Value = 258
$ws.Range("D43").Value = 2337
$ws.Range("E43").Value = 8749
$ws.Range("G43").Value = 18
$ws.Range("H43").Value = 790

# Row 47 (Dinamarca)
$ws.Range("B47").Value = 10713
$ws.Range("C47").Value = 46
$ws.Range("E47").Value = 1517

# Row 51 (Chequia)
$ws.Range("B51").Value = 8275
$ws.Range("C51").Value = 6
$ws.Range("D51").Value = 5058
$ws.Range("E51").Value = 2927
$ws.Range("F51").Value = 42
